$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "dSF" column (F) values per repulled data
$ws.Range("F3").Value = -7
$ws.Range("F5").Value = -7
$ws.Range("F6").Value = 6
$ws.Range("F8").Value = 8
$ws.Range("F9").Value = -5
$ws.Range("F10").Value = -4
$ws.Range("F13").Value = -8
$ws.Range("F14").Value = -4
$ws.Range("F17").Value = -2
$ws.Range("F18").Value = -7
